$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newValues = @("TestCheckin0001", "TestCheckin0002", "TestCheckin0003", "TestCheckin0004", "TestCheckin0005", "TestCheckin0006", "TestCheckin0007")
for ($i = 0; $i -lt 7; $i++) {
    $ws.Cells.Item($i + 2, 1).Value = $newValues[$i]
}

$ws.Columns.Item(1).ColumnWidth = 14

$ws.Range("A2:A8").Select()
